$wb = $excel.ActiveWorkbook

# --- Solids sheet: selection moves from A22 to H10, no longer the active tab ---
$wsSolids = $wb.Worksheets.Item("Solids")
$wsSolids.Range("H10").Select()

# --- Descriptions sheet: selection moves from B4 to B5 ---
$wsDescriptions = $wb.Worksheets.Item("Descriptions")
$wsDescriptions.Range("B5").Select()

# --- Other sheet: add two new rows (Reflection, Optical) and update selection ---
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("A6").Value = "Reflection"
$wsOther.Range("A7").Value = "Optical"
$wsOther.Range("A7").Select()

# "Other" becomes the active/visible tab when the workbook is reopened
$wsOther.Activate()
